# TEST EXECUTION.xlsx - "update TEST EXCETUION regist 5 #1"
# Adds two new rows (10 & 11) of regression-test data to the "EXE Regisration"
# sheet, continuing the EXE-05 test case begun in the merged D8:D9/E8:E9 block,
# mirroring the layout used for rows 8-9 (EXE-04).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clone formatting (styles/number formats/fonts/fills/borders) from the
#     existing EXE-04 block (rows 8:9) down onto the two new rows (10:11) ---
$ws.Range("D8:Q9").Copy()
$ws.Range("D10:Q11").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Rows 8/9 are tall (wrapped text); match that row height on the new rows too
$ws.Rows.Item(10).RowHeight = 31.5
$ws.Rows.Item(11).RowHeight = 31.5

# The "EXE ID"/"Scenario" columns are vertically merged across the pair of
# rows, same as D8:D9 / E8:E9
$ws.Range("D10:D11").Merge() | Out-Null
$ws.Range("E10:E11").Merge() | Out-Null

# --- Row 10: EXE-05 / EXE-TC-06 - "Register with empty email" ---
$ws.Range("D10").Value = "EXE-05"
$ws.Range("E10").Value = "Verify required fields cannot be empty"
$ws.Range("F10").Value = "EXE-TC-06"
$ws.Range("G10").Value = "Registration Module"
$ws.Range("H10").Value = "EXE-04"
$ws.Range("I10").Value = "Register with empty email"
$ws.Range("J10").Value = "Error message appears"
$ws.Range("K10").Value = "Error message " + [char]34 + "You must proviled your email" + [char]8221 + " displayed"
$ws.Range("L10").Value = "PASS"
$ws.Range("M10").Value = "-"
$ws.Range("N10").Value = "-"
$ws.Range("O10").Value = "Syaif (QA)"
$ws.Range("P10").Value = 46077
$ws.Range("Q10").Value = "Chrome v145 /`nWindows 16"

# --- Row 11: EXE-TC-07 - "Register with exactly 8 characters password" ---
$ws.Range("F11").Value = "EXE-TC-07"
$ws.Range("G11").Value = "Registration Module"
$ws.Range("H11").Value = "EXE-04"
$ws.Range("I11").Value = "Register with exactly 8 characters password"
$ws.Range("J11").Value = "Registration Successful"
$ws.Range("K11").Value = "Registration Successful"
$ws.Range("L11").Value = "PASS"
$ws.Range("M11").Value = "-"
$ws.Range("N11").Value = "-"
$ws.Range("O11").Value = "Syaif (QA)"
$ws.Range("P11").Value = 46077
$ws.Range("Q11").Value = "Chrome v145 /`nWindows 17"

# Match the saved selection/view state pointing at the new row
$ws.Range("K10").Select() | Out-Null
